$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.912.72"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "'1.908.04"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'324.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.4583"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").Value = "'0.07713"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").Value = "'1.924.74"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'5.669"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").Value = "'6.936"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "'0.07046"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'83.72"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D18").Value = "'0.000009452"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.23%  "
$ws.Range("D19").Value = "'16.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'28.898.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").Value = "'10.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").Value = "'2.096"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "'158.54"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "'19.03"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").Value = "'5.655"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").Value = "'117.64"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").Value = "'1.867"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "'0.09289"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "'0.8640"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "'5.074"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("D34").Value = "'3.090"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'0.05711"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "'1.157"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'0.02041"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("D39").Value = "'7.414"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("D41").Value = "'0.1753"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").Value = "'2.891"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.30%  "
$ws.Range("D43").Value = "'9.312"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").Value = "'0.5171"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").Value = "'2.118"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").Value = "'11.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "'0.06893"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'1.776"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("D49").Value = "'110.30"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Value = "'0.000002561"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -15.63%  "
$ws.Range("D51").Value = "'0.2857"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.35%  "
